$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
}

Set-TextValue "D2" '244.27'
Set-TextValue "E2" '-0.99%'
Set-TextValue "D3" '27.27'
Set-TextValue "E3" '3.60%'
Set-TextValue "D4" '5.111'
Set-TextValue "E4" '0.81%'
Set-TextValue "D5" '0.05659'
Set-TextValue "E5" '1.00%'
Set-TextValue "E6" '-0.39%'
Set-TextValue "D7" '0.8214'
Set-TextValue "E7" '0.98%'
Set-TextValue "D8" '0.8406'
Set-TextValue "E8" '-0.02%'
Set-TextValue "E9" '-1.38%'
Set-TextValue "D10" '0.06928'
Set-TextValue "E10" '-0.50%'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D11" '0.02989'
Set-TextValue "E11" '6.03%'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D12" '0.09395'
Set-TextValue "E12" '-0.03%'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D13" '0.001524'
Set-TextValue "E13" '0.88%'
$ws.Range("B14").Value = 'CoinExToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue "D14" '0.04207'
Set-TextValue "E14" '-10.17%'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue "D15" '0.0005982'
Set-TextValue "E15" '-93.94%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D16" '0.006139'
Set-TextValue "E16" '-0.66%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D17" '3.514'
Set-TextValue "E17" '-1.11%'
$ws.Range("B18").Value = 'GateToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue "D18" '3.003'
Set-TextValue "E18" '-1.39%'
$ws.Range("B19").Value = 'BTSEToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D19" '2.308'
Set-TextValue "E19" '8.95%'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue "D20" '0.3113'
Set-TextValue "E20" '-2.14%'
$ws.Range("B21").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C21").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D21" '0.03139'
Set-TextValue "E21" '0.75%'
Set-TextValue "D22" '0.1291'
Set-TextValue "E22" '-0.70%'
Set-TextValue "D23" '3.562'
Set-TextValue "E23" '-4.88%'
Set-TextValue "E24" '0.01%'
Set-TextValue "D25" '0.001224'
Set-TextValue "E25" '-1.98%'
Set-TextValue "D26" '0.004460'
Set-TextValue "D27" '0.00009804'
Set-TextValue "E27" '2.15%'
Set-TextValue "E28" '39.51%'
Set-TextValue "D40" '0.03652'
Set-TextValue "E40" '-0.25%'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
Set-TextValue "D41" '0.006032'
Set-TextValue "E41" '-2.53%'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-TextValue "D42" '0.1053'
Set-TextValue "E42" '-0.32%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-TextValue "D43" '0.002301'
Set-TextValue "E43" '-11.34%'
Set-TextValue "D44" '0.008977'
Set-TextValue "E44" '-5.40%'
Set-TextValue "D45" '0.00005306'
Set-TextValue "E45" '0.31%'
Set-TextValue "E47" '-36.81%'
Set-TextValue "D48" '0.002553'
Set-TextValue "E48" '23.84%'
Set-TextValue "D49" '0.00002101'
Set-TextValue "D50" '0.0002001'
